# Add more experiments. In registration process, extract the deformation
# field conditionally -> updated Dice/overlap metrics for liver.nii.gz and
# tumor.nii.gz results sheets, and highlighted the refreshed Mean row.
$wb = $excel.ActiveWorkbook

# ---- Sheet1: liver.nii.gz ----
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(2, 2).Value = 0.32956999999999997
$ws1.Cells.Item(2, 3).Value = 0
$ws1.Cells.Item(2, 4).Value = 0
$ws1.Cells.Item(3, 2).Value = 0.026280000000000001
$ws1.Cells.Item(3, 3).Value = 0.29969000000000001
$ws1.Cells.Item(3, 4).Value = 0.32244
$ws1.Cells.Item(4, 2).Value = 0.20105999999999999
$ws1.Cells.Item(4, 3).Value = 0.78681000000000001
$ws1.Cells.Item(4, 4).Value = 0.86217999999999995
$ws1.Cells.Item(5, 2).Value = 0.32612999999999998
$ws1.Cells.Item(5, 3).Value = 0.89859999999999995
$ws1.Cells.Item(5, 4).Value = 0.91452
$ws1.Cells.Item(6, 2).Value = 0
$ws1.Cells.Item(6, 3).Value = 0.71289999999999998
$ws1.Cells.Item(6, 4).Value = 0.78008999999999995
$ws1.Cells.Item(7, 2).Value = 0.0178
$ws1.Cells.Item(7, 3).Value = 0.77998000000000001
$ws1.Cells.Item(7, 4).Value = 0.76746999999999999
$ws1.Cells.Item(8, 2).Value = 0.33789999999999998
$ws1.Cells.Item(8, 3).Value = 0.77859999999999996
$ws1.Cells.Item(8, 4).Value = 0.83218000000000003
$ws1.Cells.Item(9, 2).Value = 0.16142999999999999
$ws1.Cells.Item(9, 3).Value = 0.78037000000000001
$ws1.Cells.Item(9, 4).Value = 0.82591000000000003
$ws1.Cells.Item(10, 2).Value = 0.38966000000000001
$ws1.Cells.Item(10, 3).Value = 0.52981
$ws1.Cells.Item(10, 4).Value = 0.53835999999999995
$ws1.Cells.Item(11, 2).Value = 0.098199999999999996
$ws1.Cells.Item(11, 3).Value = 0.79130999999999996
$ws1.Cells.Item(11, 4).Value = 0.83921000000000001
$ws1.Cells.Item(12, 2).Value = 0
$ws1.Cells.Item(12, 3).Value = 0
$ws1.Cells.Item(12, 4).Value = 0
$ws1.Cells.Item(13, 2).Value = 0.38966000000000001
$ws1.Cells.Item(13, 3).Value = 0.89859999999999995
$ws1.Cells.Item(13, 4).Value = 0.91452
$ws1.Cells.Item(14, 2).Value = 0.18980749999999999
$ws1.Cells.Item(14, 3).Value = 0.60472250000000005
$ws1.Cells.Item(14, 4).Value = 0.63307333333333327
$ws1.Cells.Item(15, 2).Value = 0.18980749999999999
$ws1.Cells.Item(15, 3).Value = 0.77859999999999996
$ws1.Cells.Item(15, 4).Value = 0.78008999999999995

# Highlight the refreshed Mean row C/D values (bold + yellow fill)
$meanRange1 = $ws1.Range("C14:D14")
$meanRange1.Font.Bold = $true
$meanRange1.Interior.Color = 65535

# ---- Sheet2: tumor.nii.gz ----
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(2, 2).Value = 0.13816000000000001
$ws2.Cells.Item(2, 3).Value = 0
$ws2.Cells.Item(2, 4).Value = 0
$ws2.Cells.Item(3, 2).Value = 0
$ws2.Cells.Item(3, 3).Value = 0.17029
$ws2.Cells.Item(3, 4).Value = 0.19786000000000001
$ws2.Cells.Item(4, 2).Value = 0
$ws2.Cells.Item(4, 3).Value = 0.57123999999999997
$ws2.Cells.Item(4, 4).Value = 0.82926
$ws2.Cells.Item(5, 2).Value = 0.24221999999999999
$ws2.Cells.Item(5, 3).Value = 0.85987000000000002
$ws2.Cells.Item(5, 4).Value = 0.88153000000000004
$ws2.Cells.Item(6, 2).Value = 0
$ws2.Cells.Item(6, 3).Value = 0.69177
$ws2.Cells.Item(6, 4).Value = 0.75233000000000005
$ws2.Cells.Item(7, 2).Value = 0
$ws2.Cells.Item(7, 3).Value = 0.061699999999999998
$ws2.Cells.Item(7, 4).Value = 0.062799999999999995
$ws2.Cells.Item(8, 2).Value = 0
$ws2.Cells.Item(8, 3).Value = 0.57547000000000004
$ws2.Cells.Item(8, 4).Value = 0.75199000000000005
$ws2.Cells.Item(9, 2).Value = 0.10535
$ws2.Cells.Item(9, 3).Value = 0.62211000000000005
$ws2.Cells.Item(9, 4).Value = 0.68130000000000002
$ws2.Cells.Item(10, 2).Value = 0.25814999999999999
$ws2.Cells.Item(10, 3).Value = 0.32723999999999998
$ws2.Cells.Item(10, 4).Value = 0.34838999999999998
$ws2.Cells.Item(11, 2).Value = 0
$ws2.Cells.Item(11, 3).Value = 0.67778000000000005
$ws2.Cells.Item(11, 4).Value = 0.63849999999999996
$ws2.Cells.Item(12, 2).Value = 0
$ws2.Cells.Item(12, 3).Value = 0
$ws2.Cells.Item(12, 4).Value = 0
$ws2.Cells.Item(13, 2).Value = 0.25814999999999999
$ws2.Cells.Item(13, 3).Value = 0.85987000000000002
$ws2.Cells.Item(13, 4).Value = 0.88153000000000004
$ws2.Cells.Item(14, 2).Value = 0.083502499999999993
$ws2.Cells.Item(14, 3).Value = 0.45144499999999999
$ws2.Cells.Item(14, 4).Value = 0.50212416666666659
$ws2.Cells.Item(15, 2).Value = 0
$ws2.Cells.Item(15, 3).Value = 0.57123999999999997
$ws2.Cells.Item(15, 4).Value = 0.63849999999999996

# Restore print orientation and the last-active selection on sheet1
$ws1.PageSetup.Orientation = 1
$ws1.Range("D14").Select() | Out-Null
